# "25th April 2nd update"
# Insert a new date column (25/04/2020) right after column AY (25/03/2020),
# shifting the existing 26/03/2020..31/03/2020 columns (old AZ..BF) one
# column to the right (new BA..BG), and populate the new column with the
# day's case counts. Also corrects a handful of 24/04/2020 (column AX)
# figures that were revised in the same update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new column. Everything from AZ through BF shifts right to
#    BA through BG; the freshly inserted AZ column starts out blank.
$ws.Columns("AZ").Insert()

# 2) New data for 25/04/2020 (column AZ), one entry per state row that had
#    a non-zero count that day.
$ws.Range("AZ4").Value = 61
$ws.Range("AZ9").Value = 1
$ws.Range("AZ17").Value = 15
$ws.Range("AZ28").Value = 25
$ws.Range("AZ34").Value = 57

# 3) Revised counts for 24/04/2020 (column AX), unrelated to the new
#    column but part of the same update.
$ws.Range("AX2").Value = ""
$ws.Range("AX7").Value = 53
$ws.Range("AX10").Value = 138
$ws.Range("AX16").Value = 6
$ws.Range("AX28").Value = 70
$ws.Range("AX30").Value = 13
